# feat: add 2022-Q1 data
#
# The old "总计" (totals) sheet becomes the new "2022-Q1" per-fund sheet
# (reusing its sheetId/rId and sheet-level formatting as-is), and a
# brand-new "总计" sheet is appended at the end with the refreshed
# totals table (2022-Q1 row added on top, older rows shifted down).
# The new "总计" sheet is produced by cloning an existing data sheet
# (same page setup / margins as the rest of the workbook) and then
# overwriting its cell contents, rather than Worksheets.Add() (which
# would stamp it with generic blank-workbook defaults).

$wb = $excel.ActiveWorkbook

# Donor cells already carrying the shared "header / index column" style
# (bold font + boxed border) and the plain/default (no explicit style)
# look used throughout the workbook - reused below so freshly written
# cells pick up the right look without minting new style records.
$styledDonor = $wb.Worksheets.Item("2021-Q4").Range("B1")
$plainDonor  = $wb.Worksheets.Item("2021-Q4").Range("B2")

function Set-TextValue {
    param($range, [string]$text)
    # Force the cell to store TEXT (not a number) even when the string
    # looks numeric, then strip the now-unneeded "@" number format back
    # off by pasting plain formatting over it so no stray style lingers.
    $range.NumberFormat = "@"
    $range.Value = $text
    $plainDonor.Copy()
    $range.PasteSpecial(-4122)
}

function Set-HeaderStyle {
    param($range)
    $styledDonor.Copy()
    $range.PasteSpecial(-4122)
}

# ---------------------------------------------------------------------
# 1) Repurpose the existing "总计" sheet as the new "2022-Q1" fund sheet.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

# Header row.
Set-HeaderStyle $q1.Range("B1")
$q1.Range("B1").Value = "基金代码"
Set-HeaderStyle $q1.Range("C1")
$q1.Range("C1").Value = "基金名称"
Set-HeaderStyle $q1.Range("D1")
$q1.Range("D1").Value = "基金规模"
Set-HeaderStyle $q1.Range("E1")
$q1.Range("E1").Value = "股票总仓位"
Set-HeaderStyle $q1.Range("F1")
$q1.Range("F1").Value = "仓位占比"
Set-HeaderStyle $q1.Range("G1")
$q1.Range("G1").Value = "持有市值(亿元)"
Set-HeaderStyle $q1.Range("H1")
$q1.Range("H1").Value = "仓位排名"

# Row 2 - 002423
Set-HeaderStyle $q1.Range("A2")
$q1.Range("A2").Value = 0
Set-TextValue $q1.Range("B2") "002423"
Set-TextValue $q1.Range("C2") "华宝兴业标普美国消费(QDII-LOF)美元"
Set-TextValue $q1.Range("D2") "3.62"
Set-TextValue $q1.Range("E2") "94.37"
Set-TextValue $q1.Range("F2") "3.11"
Set-TextValue $q1.Range("G2") "0.1126"
$q1.Range("H2").Value = 6

# Row 3 - 162415
Set-HeaderStyle $q1.Range("A3")
$q1.Range("A3").Value = 1
Set-TextValue $q1.Range("B3") "162415"
Set-TextValue $q1.Range("C3") "华宝标普美国消费(QDII-LOF)人民币A"
Set-TextValue $q1.Range("D3") "3.62"
Set-TextValue $q1.Range("E3") "94.37"
Set-TextValue $q1.Range("F3") "3.11"
Set-TextValue $q1.Range("G3") "0.1126"
$q1.Range("H3").Value = 6

# Row 4 - 009975
Set-HeaderStyle $q1.Range("A4")
$q1.Range("A4").Value = 2
Set-TextValue $q1.Range("B4") "009975"
Set-TextValue $q1.Range("C4") "华宝标普美国消费(QDII-LOF)人民币C"
Set-TextValue $q1.Range("D4") "0.61"
Set-TextValue $q1.Range("E4") "94.37"
Set-TextValue $q1.Range("F4") "3.11"
Set-TextValue $q1.Range("G4") "0.0190"
$q1.Range("H4").Value = 6

# Drop the now-stale rows 5:6 the old totals table left behind.
$q1.Range("A5:H6").Clear()

# ---------------------------------------------------------------------
# 2) Append a brand-new "总计" sheet with the refreshed totals table.
#    Clone an existing sheet (matching page setup/margins) rather than
#    Worksheets.Add(), then overwrite its cells.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q4").Copy([System.Reflection.Missing]::Value, $q1)
$total = $wb.Worksheets.Item($q1.Index + 1)
$total.Name = "总计"

# Drop the template's extra E:H columns - the totals sheet only needs A:D.
$total.Range("E1:H4").Clear()

Set-HeaderStyle $total.Range("B1")
$total.Range("B1").Value = "日期"
Set-HeaderStyle $total.Range("C1")
$total.Range("C1").Value = "持有数量(只)"
Set-HeaderStyle $total.Range("D1")
$total.Range("D1").Value = "持有市值(亿元)"

Set-HeaderStyle $total.Range("A2")
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.24

Set-HeaderStyle $total.Range("A3")
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 3
$total.Range("D3").Value = 0.34

Set-HeaderStyle $total.Range("A4")
$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 3
$total.Range("D4").Value = 0.25

Set-HeaderStyle $total.Range("A5")
$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q2"
$total.Range("C5").Value = 3
$total.Range("D5").Value = 0.22

Set-HeaderStyle $total.Range("A6")
$total.Range("A6").Value = 4
$total.Range("B6").Value = "2021-Q1"
$total.Range("C6").Value = 5
$total.Range("D6").Value = 0.23

Set-HeaderStyle $total.Range("A7")
$total.Range("A7").Value = 5
$total.Range("B7").Value = "2020-Q4"
$total.Range("C7").Value = 3
$total.Range("D7").Value = 0.17

# Restore the originally-active tab so tabSelected stays on "2020-Q4".
$wb.Worksheets.Item("2020-Q4").Activate()
